$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the existing row 172 ("Carson" block),
# pushing old rows 172:263 down to 174:265 (matches the new dimension A1:T265).
$ws.Rows.Item(172).Resize(2).Insert()

# --- New row 172: Fruta, Feria Lagunitas de Puerto Montt - Durazno / Carson / Especial ---
$ws.Cells.Item(172, 1).Value = 4
$ws.Cells.Item(172, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(172, 3).Value = "Los Lagos"
$ws.Cells.Item(172, 4).Value = 44960
$ws.Cells.Item(172, 5).Value = 10
$ws.Cells.Item(172, 6).Value = "Fruta"
$ws.Cells.Item(172, 7).Value = 100103
$ws.Cells.Item(172, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(172, 9).Value = 100103004
$ws.Cells.Item(172, 10).Value = "Durazno"
$ws.Cells.Item(172, 11).Value = "Carson"
$ws.Cells.Item(172, 12).Value = "Especial"
$ws.Cells.Item(172, 13).Value = 200
$ws.Cells.Item(172, 14).Value = 22000
$ws.Cells.Item(172, 15).Value = 22000
$ws.Cells.Item(172, 16).Value = 22000
$ws.Cells.Item(172, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(172, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(172, 19).Value = 1571
$ws.Cells.Item(172, 20).Value = 14

# --- New row 173: Fruta, Feria Lagunitas de Puerto Montt - Durazno / Carson / Primera ---
$ws.Cells.Item(173, 1).Value = 4
$ws.Cells.Item(173, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(173, 3).Value = "Los Lagos"
$ws.Cells.Item(173, 4).Value = 44960
$ws.Cells.Item(173, 5).Value = 10
$ws.Cells.Item(173, 6).Value = "Fruta"
$ws.Cells.Item(173, 7).Value = 100103
$ws.Cells.Item(173, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(173, 9).Value = 100103004
$ws.Cells.Item(173, 10).Value = "Durazno"
$ws.Cells.Item(173, 11).Value = "Carson"
$ws.Cells.Item(173, 12).Value = "Primera"
$ws.Cells.Item(173, 13).Value = 400
$ws.Cells.Item(173, 14).Value = 18000
$ws.Cells.Item(173, 15).Value = 19000
$ws.Cells.Item(173, 16).Value = 18500
$ws.Cells.Item(173, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(173, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(173, 19).Value = 1321
$ws.Cells.Item(173, 20).Value = 14
